# Updated cryptos list on Sat Dec 30 10:53:23 UTC 2023 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values,
# plus a coin-list change in row 51 (Cronos -> TrustWalletToken).
#
# NOTE: Price values are stored as plain text in the sheet (e.g. "315.60"),
# not as numbers, so we prefix them with a leading apostrophe to force
# Excel to keep them as literal text instead of auto-converting them to
# numeric values (which would silently drop meaningful trailing zeros,
# e.g. "315.60" -> 315.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.983.80"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "'2.288.49"
$ws.Range("E3").Value = "  -3.39%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'315.60"
$ws.Range("E5").Value = "  -0.88%  "

$ws.Range("D6").Value = "'102.69"
$ws.Range("E6").Value = "  -5.66%  "

$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  -1.46%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  -3.62%  "

$ws.Range("D10").Value = "'38.86"
$ws.Range("E10").Value = "  -7.38%  "

$ws.Range("D11").Value = "'0.0904"
$ws.Range("E11").Value = "  -2.69%  "

$ws.Range("D12").Value = "'8.23"
$ws.Range("E12").Value = "  -3.94%  "

$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").Value = "'0.960"
$ws.Range("E14").Value = "  -4.80%  "

$ws.Range("D15").Value = "'15.25"
$ws.Range("E15").Value = "  -5.62%  "

$ws.Range("D16").Value = "'2.633.63"
$ws.Range("E16").Value = "  -3.42%  "

$ws.Range("D17").Value = "'2.277.86"
$ws.Range("E17").Value = "  -3.19%  "

$ws.Range("D18").Value = "'41.854.61"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").Value = "'7.43"
$ws.Range("E19").Value = "  -3.37%  "

$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").Value = "'3.63"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").Value = "'73.13"
$ws.Range("E22").Value = "  -4.09%  "

$ws.Range("D23").Value = "'280.39"
$ws.Range("E23").Value = "  +8.97%  "

$ws.Range("D24").Value = "'10.15"
$ws.Range("E24").Value = "  +7.44%  "

$ws.Range("E25").Value = "  -4.35%  "

$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").Value = "'2.40"
$ws.Range("E27").Value = "  +6.75%  "

$ws.Range("D28").Value = "'10.71"
$ws.Range("E28").Value = "  -6.64%  "

$ws.Range("D29").Value = "'22.86"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("D30").Value = "'163.40"
$ws.Range("E30").Value = "  -5.14%  "

$ws.Range("D31").Value = "'35.18"
$ws.Range("E31").Value = "  -4.97%  "

$ws.Range("D32").Value = "'0.0869"
$ws.Range("E32").Value = "  -3.29%  "

$ws.Range("D33").Value = "'2.85"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("D34").Value = "'5.80"
$ws.Range("E34").Value = "  -4.42%  "

$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("E36").Value = "  -4.47%  "

$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  -3.22%  "

$ws.Range("E38").Value = "  +7.18%  "

$ws.Range("E39").Value = "  -5.42%  "

$ws.Range("D40").Value = "'3.63"
$ws.Range("E40").Value = "  -8.02%  "

$ws.Range("D41").Value = "'99.82"
$ws.Range("E41").Value = "  +13.24%  "

$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("D43").Value = "'69.34"
$ws.Range("E43").Value = "  -3.54%  "

$ws.Range("D45").Value = "'0.224"
$ws.Range("E45").Value = "  -7.81%  "

$ws.Range("D46").Value = "'115.25"
$ws.Range("E46").Value = "  +1.64%  "

$ws.Range("D47").Value = "'11.83"
$ws.Range("E47").Value = "  -4.31%  "

$ws.Range("D48").Value = "'8.91"
$ws.Range("E48").Value = "  -3.68%  "

$ws.Range("D49").Value = "'76.03"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'5.24"
$ws.Range("E50").Value = "  -6.33%  "

# Row 51: coin changed from Cronos to TrustWalletToken
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.25"
$ws.Range("E51").Value = "  -4.47%  "
